$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 10) to the customers data, mirroring an existing
# "good" record but using a multi-space string for the name field to
# exercise the "empty" (blank/whitespace) unauthorized-field test case.

$ws.Range("A10").Value = "       "
$ws.Range("B10").Value = 56
$ws.Range("C10").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = "Wine"
$ws.Range("E10").Value = "13/04/2001"
$ws.Range("F10").Value = 1.9

[void]$ws.Range("C11").Select()
